# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures to the Goblin Profits workbook
# (columns H-N on various rows across all 8 job sheets), per the commit diff.

$wb = $excel.ActiveWorkbook

# ALC row 5: Met a Sticky End / Animal Glue
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 109
$ws.Range("I5").Value = 116.666664
$ws.Range("K5").Value = 116.666664
$ws.Range("M5").Value = -1.666663999999997

# ALC row 18: You Grow, Girl / Growth Formula Beta
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 7325.5
$ws.Range("J18").Value = 3480.8
$ws.Range("L18").Value = 3480.8
$ws.Range("N18").Value = -4048.8

# ALC row 40: Stuck in the Moment / Horn Glue
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3536.077
$ws.Range("I40").Value = 1269.25
$ws.Range("J40").Value = 4543.5557
$ws.Range("K40").Value = 1269.25
$ws.Range("L40").Value = 4543.5557
$ws.Range("M40").Value = -1094.25
$ws.Range("N40").Value = -4893.5557

# ALC row 43: Growing Is Knowing / Growth Formula Gamma
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5072.278
$ws.Range("I43").Value = 2361
$ws.Range("J43").Value = 5846.9287
$ws.Range("K43").Value = 2361
$ws.Range("L43").Value = 5846.9287
$ws.Range("M43").Value = -2292
$ws.Range("N43").Value = -5984.9287

# ALC row 98: The Dotted Line / Enchanted Durium Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1870.4445
$ws.Range("I98").Value = 1521.2084
$ws.Range("K98").Value = 1521.2084
$ws.Range("M98").Value = -23.20839999999998

# ALC row 106: Making Your Mark / Enchanted Palladium Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3100.75
$ws.Range("I106").Value = 1201.5
$ws.Range("K106").Value = 1201.5
$ws.Range("M106").Value = -570.5

# ALC row 113: Amaro Kart / Starch Glue
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 6044.3794
$ws.Range("J113").Value = 3748.1667
$ws.Range("L113").Value = 3748.1667
$ws.Range("N113").Value = -10256.1667

# ALC row 122: Wishful Inking / Enchanted High Durium Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1870.4445
$ws.Range("I122").Value = 1521.2084
$ws.Range("K122").Value = 4563.6252
$ws.Range("M122").Value = -2113.6252

# ALC row 129: Practical Command / Commanding Craftsman's Draught
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1194.4166
$ws.Range("J129").Value = 2184.889
$ws.Range("L129").Value = 6554.667
$ws.Range("N129").Value = -16554.667

# ARM row 32: Ingot We Trust / Steel Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2835.125
$ws.Range("I32").Value = 3088.2927
$ws.Range("K32").Value = 3088.2927
$ws.Range("M32").Value = -2801.2927

# ARM row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4492.381
$ws.Range("I102").Value = 2289.3333
$ws.Range("K102").Value = 2289.3333
$ws.Range("M102").Value = -667.3332999999998

# ARM row 103: Sweeping the Legs / Doman Steel Greaves of Striking
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# BSM row 86: Through Thick and Thin / Adamantite Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3135.7
$ws.Range("I86").Value = 3091.1428
$ws.Range("K86").Value = 3091.1428
$ws.Range("M86").Value = -1968.1428

# BSM row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3135.7
$ws.Range("I89").Value = 3091.1428
$ws.Range("K89").Value = 15455.714
$ws.Range("M89").Value = -9839.714

# BSM row 94: High Steal / High Steel Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2120.9333
$ws.Range("I94").Value = 2648.762
$ws.Range("J94").Value = 889.3333
$ws.Range("K94").Value = 2648.762
$ws.Range("L94").Value = 889.3333
$ws.Range("M94").Value = -2197.762
$ws.Range("N94").Value = -1791.3333

# BSM row 134: Ruthenium Supremium / Ruthenium Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3951.7896
$ws.Range("I134").Value = 4252.769
$ws.Range("K134").Value = 12758.307
$ws.Range("M134").Value = -10223.307

# CRP row 6: Got Your Back / Square Maple Shield
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1000.25
$ws.Range("I6").Value = 1000.25
$ws.Range("K6").Value = 1000.25
$ws.Range("M6").Value = -887.25

# CRP row 17: Say It with Spears / Feathered Harpoon
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 9999.5
$ws.Range("I17").Value = 9999
$ws.Range("K17").Value = 9999
$ws.Range("M17").Value = -9825

# CRP row 22: Driving Up the Wall / Elm Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1318.84
$ws.Range("I22").Value = 884
$ws.Range("K22").Value = 884
$ws.Range("M22").Value = -534

# CRP row 31: Wall Not Found / Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3017.1428
$ws.Range("I31").Value = 1457.375
$ws.Range("J31").Value = 6420.273
$ws.Range("K31").Value = 1457.375
$ws.Range("L31").Value = 6420.273
$ws.Range("M31").Value = -1162.375
$ws.Range("N31").Value = -7010.273

# CRP row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3017.1428
$ws.Range("I34").Value = 1457.375
$ws.Range("J34").Value = 6420.273
$ws.Range("K34").Value = 1457.375
$ws.Range("L34").Value = 6420.273
$ws.Range("M34").Value = -1255.375
$ws.Range("N34").Value = -6824.273

# CRP row 50: The Arsenal of Theocracy / Cobalt Halberd
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 54997.5
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 54997.5
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 54997.5
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -56247.5

# CRP row 55: Ready for a Rematch / Mythril Lance
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 8899.666999999999
$ws.Range("I55").Value = 8899.666999999999
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 8899.666999999999
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -8584.666999999999
$ws.Range("N55").ClearContents()

# CRP row 60: Bowing to Greater Power / Yew Longbow
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 26332
$ws.Range("I60").Value = 13857
$ws.Range("K60").Value = 13857
$ws.Range("M60").Value = -13346

# CRP row 62: Splinter in the Sewers / Cedar Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 7549.6665
$ws.Range("I62").Value = 7559.6
$ws.Range("K62").Value = 7559.6
$ws.Range("M62").Value = -6935.6

# CRP row 64: Almost as Fun as Slingshotting Birds / Cedar Longbow
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496

# CRP row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 7549.6665
$ws.Range("I65").Value = 7559.6
$ws.Range("K65").Value = 37798
$ws.Range("M65").Value = -34678

# CRP row 67: Living Bow to Mouth (L) / Cedar Longbow
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716

# CRP row 74: License to Heal / Dark Chestnut Rod
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# CRP row 77: Purified Polyrhythm (L) / Dark Chestnut Rod
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# CRP row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2067.9707
$ws.Range("I122").Value = 1752.9584
$ws.Range("J122").Value = 2824
$ws.Range("K122").Value = 5258.8752
$ws.Range("L122").Value = 8472
$ws.Range("M122").Value = -2808.8752
$ws.Range("N122").Value = -13372

# CUL row 131: The Mountain Steeped / Tsai tou Vounou
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2383.12
$ws.Range("I131").Value = 1075.6
$ws.Range("J131").Value = 2710
$ws.Range("K131").Value = 3226.8
$ws.Range("L131").Value = 8130
$ws.Range("M131").Value = 1813.2
$ws.Range("N131").Value = -18210

# CUL row 140: Sweet, Sweet Bean Juice / Mesquite Juice
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2594.7144
$ws.Range("I140").Value = 2360.5
$ws.Range("K140").Value = 7081.5
$ws.Range("M140").Value = -1901.5

# GSM row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4947.9375
$ws.Range("I97").Value = 1865.1072
$ws.Range("J97").Value = 26527.75
$ws.Range("K97").Value = 1865.1072
$ws.Range("L97").Value = 26527.75
$ws.Range("M97").Value = -1369.1072
$ws.Range("N97").Value = -27519.75

# GSM row 132: On Board for Lar / Lar Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2632.4167
$ws.Range("J132").Value = 3399.6667
$ws.Range("L132").Value = 10199.0001
$ws.Range("N132").Value = -15259.0001

# LTW row 16: Saddle Sore / Hard Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3859.8948
$ws.Range("I16").Value = 3453
$ws.Range("J16").Value = 4999.2
$ws.Range("K16").Value = 3453
$ws.Range("L16").Value = 4999.2
$ws.Range("M16").Value = -3283
$ws.Range("N16").Value = -5339.2

# LTW row 46: Supply Side Logic / Boar Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2357.2122
$ws.Range("I46").Value = 1289.9231
$ws.Range("J46").Value = 3050.95
$ws.Range("K46").Value = 1289.9231
$ws.Range("L46").Value = 3050.95
$ws.Range("M46").Value = -1101.9231
$ws.Range("N46").Value = -3426.95

# LTW row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 7213.963
$ws.Range("J68").Value = 8388.777
$ws.Range("L68").Value = 8388.777
$ws.Range("N68").Value = -9886.777

# LTW row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 7213.963
$ws.Range("J71").Value = 8388.777
$ws.Range("L71").Value = 41943.885
$ws.Range("N71").Value = -49431.885

# LTW row 93: Hide to Go Seek / Gagana Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2364.3948
$ws.Range("I93").Value = 1664
$ws.Range("K93").Value = 1664
$ws.Range("M93").Value = -416

# WVR row 74: Clothing the Naked Truth / Ramie Robe of Casting
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 31712.25
$ws.Range("J74").Value = 28616.334
$ws.Range("L74").Value = 28616.334
$ws.Range("N74").Value = -30488.334

# WVR row 77: When in Robes (L) / Ramie Robe of Casting
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H77").Value = 31712.25
$ws.Range("J77").Value = 28616.334
$ws.Range("L77").Value = 85849.00199999999
$ws.Range("N77").Value = -95209.00199999999

# WVR row 100: Of Great Import / Kudzu Thread
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1813.6666
$ws.Range("I100").Value = 1439.8334
$ws.Range("K100").Value = 2879.6668
$ws.Range("M100").Value = -2338.6668

# WVR row 126: A Polished Purchase / Snow Linen
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1682.3226
$ws.Range("I126").Value = 1258.9584
$ws.Range("K126").Value = 3776.8752
$ws.Range("M126").Value = -1306.8752

# WVR row 136: Weaving the Envelope / Sarcenet Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 967.0732
$ws.Range("I136").Value = 909.97144
$ws.Range("J136").Value = 1300.1666
$ws.Range("K136").Value = 2729.91432
$ws.Range("L136").Value = 3900.4998
$ws.Range("M136").Value = -179.9143199999999
$ws.Range("N136").Value = -9000.4998

